{"js": "// Add a 1x2 \"saved-by\" / author table at the end of the document body,\n// followed by an empty paragraph (mirrors python-docx's\n// Document.add_table() + iter_inner_content() usage described in the\n// commit message), even though this document has no sections.\n\nconst body = context.document.body;\n\n// Minimal OOXML package containing just the table + trailing empty\n// paragraph we want appended. Word (and this shim) splices the body\n// content of the supplied part into the target location.\nconst ooxml =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>\\n' +\n  '<?mso-application progid=\"Word.Document\"?>\\n' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:tbl>' +\n  '<w:tblGrid>' +\n  '<w:gridCol w:w=\"4788\"/>' +\n  '<w:gridCol w:w=\"4788\"/>' +\n  '</w:tblGrid>' +\n  '<w:tr>' +\n  '<w:tc><w:p><w:r><w:t>saved-by</w:t></w:r></w:p></w:tc>' +\n  '<w:tc><w:p><w:r><w:t>Dennis Forsythe</w:t></w:r></w:p></w:tc>' +\n  '</w:tr>' +\n  '</w:tbl>' +\n  '<w:p/>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nbody.insertOoxml(ooxml, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Add a 1x2 \"saved-by\" / author table at the end of the document body,\n# followed by an empty paragraph (mirrors python-docx's\n# Document.add_table() + iter_inner_content() usage described in the\n# commit message), even though this document has no sections.\n\n$d = $word.ActiveDocument\n\n# Collapsed range at the very end of the document's main story.\n$r = $d.Range($d.Content.End, $d.Content.End)\n\n# Minimal WordOpenXML package containing just the table + trailing\n# empty paragraph we want appended. InsertXML replaces the (collapsed)\n# range's contents with the supplied markup, effectively inserting it.\n$ooxml = '<?xml version=\"1.0\" standalone=\"yes\"?><?mso-application progid=\"Word.Document\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:tbl><w:tblGrid><w:gridCol w:w=\"4788\"/><w:gridCol w:w=\"4788\"/></w:tblGrid><w:tr><w:tc><w:p><w:r><w:t>saved-by</w:t></w:r></w:p></w:tc><w:tc><w:p><w:r><w:t>Dennis Forsythe</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$r.InsertXML($ooxml)\n"}
